# Add data for 2022-03-22: update sheet title, March label, and figures for
# March row (row 4) and Total row (row 5).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet (tab title) to reflect the new "through" date.
$ws.Name = "Through 2022-03-14"

# Update the label cell for the March row.
$ws.Range("A4").Value = "March (through 03-14)"

# Update March row (row 4) figures.
$ws.Cells.Item(4, 2).Value = 12   # B4
$ws.Cells.Item(4, 5).Value = 26   # E4
$ws.Cells.Item(4, 6).Value = 14   # F4
$ws.Cells.Item(4, 8).Value = 38   # H4
$ws.Cells.Item(4, 9).Value = 65   # I4

# Update Total row (row 5) figures.
$ws.Cells.Item(5, 2).Value = 49    # B5
$ws.Cells.Item(5, 5).Value = 163   # E5
$ws.Cells.Item(5, 6).Value = 93    # F5
$ws.Cells.Item(5, 8).Value = 380   # H5
$ws.Cells.Item(5, 9).Value = 365   # I5
